$wb = $excel.ActiveWorkbook

# --- PIR sheet: add rows 615-626 ---
$wsPIR = $wb.Worksheets.Item(2)
$wsPIR.Cells.Item(615, 1).NumberFormat = "@"
$wsPIR.Cells.Item(615, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(615, 2).Value = "10:32:24"
$wsPIR.Cells.Item(615, 3).Value = "10:00"
$wsPIR.Cells.Item(615, 4).Value = "Bathroom"
$wsPIR.Cells.Item(615, 5).Value = "No Motion"
$wsPIR.Cells.Item(615, 6).Value = "Inactive"

$wsPIR.Cells.Item(616, 1).NumberFormat = "@"
$wsPIR.Cells.Item(616, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(616, 2).Value = "10:32:30"
$wsPIR.Cells.Item(616, 3).Value = "10:00"
$wsPIR.Cells.Item(616, 4).Value = "Bathroom"
$wsPIR.Cells.Item(616, 5).Value = "No Motion"
$wsPIR.Cells.Item(616, 6).Value = "Inactive"

$wsPIR.Cells.Item(617, 1).NumberFormat = "@"
$wsPIR.Cells.Item(617, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(617, 2).Value = "10:32:35"
$wsPIR.Cells.Item(617, 3).Value = "10:00"
$wsPIR.Cells.Item(617, 4).Value = "Bathroom"
$wsPIR.Cells.Item(617, 5).Value = "No Motion"
$wsPIR.Cells.Item(617, 6).Value = "Inactive"

$wsPIR.Cells.Item(618, 1).NumberFormat = "@"
$wsPIR.Cells.Item(618, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(618, 2).Value = "10:32:40"
$wsPIR.Cells.Item(618, 3).Value = "10:00"
$wsPIR.Cells.Item(618, 4).Value = "Bathroom"
$wsPIR.Cells.Item(618, 5).Value = "No Motion"
$wsPIR.Cells.Item(618, 6).Value = "Inactive"

$wsPIR.Cells.Item(619, 1).NumberFormat = "@"
$wsPIR.Cells.Item(619, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(619, 2).Value = "10:32:45"
$wsPIR.Cells.Item(619, 3).Value = "10:00"
$wsPIR.Cells.Item(619, 4).Value = "Bathroom"
$wsPIR.Cells.Item(619, 5).Value = "No Motion"
$wsPIR.Cells.Item(619, 6).Value = "Inactive"

$wsPIR.Cells.Item(620, 1).NumberFormat = "@"
$wsPIR.Cells.Item(620, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(620, 2).Value = "10:32:50"
$wsPIR.Cells.Item(620, 3).Value = "10:00"
$wsPIR.Cells.Item(620, 4).Value = "Bathroom"
$wsPIR.Cells.Item(620, 5).Value = "No Motion"
$wsPIR.Cells.Item(620, 6).Value = "Inactive"

$wsPIR.Cells.Item(621, 1).NumberFormat = "@"
$wsPIR.Cells.Item(621, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(621, 2).Value = "10:32:55"
$wsPIR.Cells.Item(621, 3).Value = "10:00"
$wsPIR.Cells.Item(621, 4).Value = "Bathroom"
$wsPIR.Cells.Item(621, 5).Value = "No Motion"
$wsPIR.Cells.Item(621, 6).Value = "Inactive"

$wsPIR.Cells.Item(622, 1).NumberFormat = "@"
$wsPIR.Cells.Item(622, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(622, 2).Value = "10:33:00"
$wsPIR.Cells.Item(622, 3).Value = "10:00"
$wsPIR.Cells.Item(622, 4).Value = "Bathroom"
$wsPIR.Cells.Item(622, 5).Value = "No Motion"
$wsPIR.Cells.Item(622, 6).Value = "Inactive"

$wsPIR.Cells.Item(623, 1).NumberFormat = "@"
$wsPIR.Cells.Item(623, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(623, 2).Value = "10:33:05"
$wsPIR.Cells.Item(623, 3).Value = "10:00"
$wsPIR.Cells.Item(623, 4).Value = "Bathroom"
$wsPIR.Cells.Item(623, 5).Value = "No Motion"
$wsPIR.Cells.Item(623, 6).Value = "Inactive"

$wsPIR.Cells.Item(624, 1).NumberFormat = "@"
$wsPIR.Cells.Item(624, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(624, 2).Value = "10:33:10"
$wsPIR.Cells.Item(624, 3).Value = "10:00"
$wsPIR.Cells.Item(624, 4).Value = "Bathroom"
$wsPIR.Cells.Item(624, 5).Value = "No Motion"
$wsPIR.Cells.Item(624, 6).Value = "Inactive"

$wsPIR.Cells.Item(625, 1).NumberFormat = "@"
$wsPIR.Cells.Item(625, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(625, 2).Value = "10:33:15"
$wsPIR.Cells.Item(625, 3).Value = "10:00"
$wsPIR.Cells.Item(625, 4).Value = "Bathroom"
$wsPIR.Cells.Item(625, 5).Value = "No Motion"
$wsPIR.Cells.Item(625, 6).Value = "Inactive"

$wsPIR.Cells.Item(626, 1).NumberFormat = "@"
$wsPIR.Cells.Item(626, 1).Value = "2026-02-06"
$wsPIR.Cells.Item(626, 2).Value = "10:33:20"
$wsPIR.Cells.Item(626, 3).Value = "10:00"
$wsPIR.Cells.Item(626, 4).Value = "Bathroom"
$wsPIR.Cells.Item(626, 5).Value = "No Motion"
$wsPIR.Cells.Item(626, 6).Value = "Inactive"

# --- Humidity sheet: add rows 441-449 ---
$wsHum = $wb.Worksheets.Item(3)
$wsHum.Cells.Item(441, 1).NumberFormat = "@"
$wsHum.Cells.Item(441, 1).Value = "2026-02-06"
$wsHum.Cells.Item(441, 2).Value = "10:32:26"
$wsHum.Cells.Item(441, 3).Value = "10:00"
$wsHum.Cells.Item(441, 4).Value = "Bathroom"
$wsHum.Cells.Item(441, 5).NumberFormat = "@"
$wsHum.Cells.Item(441, 5).Value = "67.8%"
$wsHum.Cells.Item(441, 6).Value = "Active"

$wsHum.Cells.Item(442, 1).NumberFormat = "@"
$wsHum.Cells.Item(442, 1).Value = "2026-02-06"
$wsHum.Cells.Item(442, 2).Value = "10:32:31"
$wsHum.Cells.Item(442, 3).Value = "10:00"
$wsHum.Cells.Item(442, 4).Value = "Bathroom"
$wsHum.Cells.Item(442, 5).NumberFormat = "@"
$wsHum.Cells.Item(442, 5).Value = "67.6%"
$wsHum.Cells.Item(442, 6).Value = "Active"

$wsHum.Cells.Item(443, 1).NumberFormat = "@"
$wsHum.Cells.Item(443, 1).Value = "2026-02-06"
$wsHum.Cells.Item(443, 2).Value = "10:32:36"
$wsHum.Cells.Item(443, 3).Value = "10:00"
$wsHum.Cells.Item(443, 4).Value = "Bathroom"
$wsHum.Cells.Item(443, 5).NumberFormat = "@"
$wsHum.Cells.Item(443, 5).Value = "67.5%"
$wsHum.Cells.Item(443, 6).Value = "Active"

$wsHum.Cells.Item(444, 1).NumberFormat = "@"
$wsHum.Cells.Item(444, 1).Value = "2026-02-06"
$wsHum.Cells.Item(444, 2).Value = "10:32:46"
$wsHum.Cells.Item(444, 3).Value = "10:00"
$wsHum.Cells.Item(444, 4).Value = "Bathroom"
$wsHum.Cells.Item(444, 5).NumberFormat = "@"
$wsHum.Cells.Item(444, 5).Value = "67.3%"
$wsHum.Cells.Item(444, 6).Value = "Active"

$wsHum.Cells.Item(445, 1).NumberFormat = "@"
$wsHum.Cells.Item(445, 1).Value = "2026-02-06"
$wsHum.Cells.Item(445, 2).Value = "10:32:51"
$wsHum.Cells.Item(445, 3).Value = "10:00"
$wsHum.Cells.Item(445, 4).Value = "Bathroom"
$wsHum.Cells.Item(445, 5).NumberFormat = "@"
$wsHum.Cells.Item(445, 5).Value = "67.3%"
$wsHum.Cells.Item(445, 6).Value = "Active"

$wsHum.Cells.Item(446, 1).NumberFormat = "@"
$wsHum.Cells.Item(446, 1).Value = "2026-02-06"
$wsHum.Cells.Item(446, 2).Value = "10:32:56"
$wsHum.Cells.Item(446, 3).Value = "10:00"
$wsHum.Cells.Item(446, 4).Value = "Bathroom"
$wsHum.Cells.Item(446, 5).NumberFormat = "@"
$wsHum.Cells.Item(446, 5).Value = "67.3%"
$wsHum.Cells.Item(446, 6).Value = "Active"

$wsHum.Cells.Item(447, 1).NumberFormat = "@"
$wsHum.Cells.Item(447, 1).Value = "2026-02-06"
$wsHum.Cells.Item(447, 2).Value = "10:33:01"
$wsHum.Cells.Item(447, 3).Value = "10:00"
$wsHum.Cells.Item(447, 4).Value = "Bathroom"
$wsHum.Cells.Item(447, 5).NumberFormat = "@"
$wsHum.Cells.Item(447, 5).Value = "67.4%"
$wsHum.Cells.Item(447, 6).Value = "Active"

$wsHum.Cells.Item(448, 1).NumberFormat = "@"
$wsHum.Cells.Item(448, 1).Value = "2026-02-06"
$wsHum.Cells.Item(448, 2).Value = "10:33:16"
$wsHum.Cells.Item(448, 3).Value = "10:00"
$wsHum.Cells.Item(448, 4).Value = "Bathroom"
$wsHum.Cells.Item(448, 5).NumberFormat = "@"
$wsHum.Cells.Item(448, 5).Value = "67.2%"
$wsHum.Cells.Item(448, 6).Value = "Active"

$wsHum.Cells.Item(449, 1).NumberFormat = "@"
$wsHum.Cells.Item(449, 1).Value = "2026-02-06"
$wsHum.Cells.Item(449, 2).Value = "10:33:21"
$wsHum.Cells.Item(449, 3).Value = "10:00"
$wsHum.Cells.Item(449, 4).Value = "Bathroom"
$wsHum.Cells.Item(449, 5).NumberFormat = "@"
$wsHum.Cells.Item(449, 5).Value = "67.2%"
$wsHum.Cells.Item(449, 6).Value = "Active"

# --- Temperature sheet: add rows 440-448 ---
$wsTemp = $wb.Worksheets.Item(4)
$wsTemp.Cells.Item(440, 1).NumberFormat = "@"
$wsTemp.Cells.Item(440, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(440, 2).Value = "10:32:27"
$wsTemp.Cells.Item(440, 3).Value = "10:00"
$wsTemp.Cells.Item(440, 4).Value = "Bathroom"
$wsTemp.Cells.Item(440, 5).Value = "28.5C"
$wsTemp.Cells.Item(440, 6).Value = "Active"

$wsTemp.Cells.Item(441, 1).NumberFormat = "@"
$wsTemp.Cells.Item(441, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(441, 2).Value = "10:32:32"
$wsTemp.Cells.Item(441, 3).Value = "10:00"
$wsTemp.Cells.Item(441, 4).Value = "Bathroom"
$wsTemp.Cells.Item(441, 5).Value = "28.5C"
$wsTemp.Cells.Item(441, 6).Value = "Active"

$wsTemp.Cells.Item(442, 1).NumberFormat = "@"
$wsTemp.Cells.Item(442, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(442, 2).Value = "10:32:37"
$wsTemp.Cells.Item(442, 3).Value = "10:00"
$wsTemp.Cells.Item(442, 4).Value = "Bathroom"
$wsTemp.Cells.Item(442, 5).Value = "28.5C"
$wsTemp.Cells.Item(442, 6).Value = "Active"

$wsTemp.Cells.Item(443, 1).NumberFormat = "@"
$wsTemp.Cells.Item(443, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(443, 2).Value = "10:32:47"
$wsTemp.Cells.Item(443, 3).Value = "10:00"
$wsTemp.Cells.Item(443, 4).Value = "Bathroom"
$wsTemp.Cells.Item(443, 5).Value = "28.5C"
$wsTemp.Cells.Item(443, 6).Value = "Active"

$wsTemp.Cells.Item(444, 1).NumberFormat = "@"
$wsTemp.Cells.Item(444, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(444, 2).Value = "10:32:52"
$wsTemp.Cells.Item(444, 3).Value = "10:00"
$wsTemp.Cells.Item(444, 4).Value = "Bathroom"
$wsTemp.Cells.Item(444, 5).Value = "28.5C"
$wsTemp.Cells.Item(444, 6).Value = "Active"

$wsTemp.Cells.Item(445, 1).NumberFormat = "@"
$wsTemp.Cells.Item(445, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(445, 2).Value = "10:32:57"
$wsTemp.Cells.Item(445, 3).Value = "10:00"
$wsTemp.Cells.Item(445, 4).Value = "Bathroom"
$wsTemp.Cells.Item(445, 5).Value = "28.5C"
$wsTemp.Cells.Item(445, 6).Value = "Active"

$wsTemp.Cells.Item(446, 1).NumberFormat = "@"
$wsTemp.Cells.Item(446, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(446, 2).Value = "10:33:03"
$wsTemp.Cells.Item(446, 3).Value = "10:00"
$wsTemp.Cells.Item(446, 4).Value = "Bathroom"
$wsTemp.Cells.Item(446, 5).Value = "28.6C"
$wsTemp.Cells.Item(446, 6).Value = "Active"

$wsTemp.Cells.Item(447, 1).NumberFormat = "@"
$wsTemp.Cells.Item(447, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(447, 2).Value = "10:33:18"
$wsTemp.Cells.Item(447, 3).Value = "10:00"
$wsTemp.Cells.Item(447, 4).Value = "Bathroom"
$wsTemp.Cells.Item(447, 5).Value = "28.5C"
$wsTemp.Cells.Item(447, 6).Value = "Active"

$wsTemp.Cells.Item(448, 1).NumberFormat = "@"
$wsTemp.Cells.Item(448, 1).Value = "2026-02-06"
$wsTemp.Cells.Item(448, 2).Value = "10:33:23"
$wsTemp.Cells.Item(448, 3).Value = "10:00"
$wsTemp.Cells.Item(448, 4).Value = "Bathroom"
$wsTemp.Cells.Item(448, 5).Value = "28.6C"
$wsTemp.Cells.Item(448, 6).Value = "Active"

